# Update for 苏州-漫展信息.xlsx
#
# The edit removes the first data row (2024-04-04 "Anime LIVE" expo) from
# every sheet that listed it ("展览" and "全部类型"), shifting the rest of
# the rows up by one, renumbering the running index in column A, and
# refreshing the "想去人数" (F column) view-count numbers for the events
# that remain (a re-scrape of live stats). The "演出" and "本地生活"
# sheets are untouched because they never contained that row / are not
# touched by the diff.

$wb = $excel.ActiveWorkbook

# Map of bilibili show id (as found at the end of the Link / H column URL)
# -> refreshed "想去人数" (F column) value, taken from the target diff.
$idToNewF = @{
    82042 = 561
    80789 = 495
    81879 = 1230
    78666 = 1060
    77196 = 13982
    79789 = 15207
    83576 = 6
    83575 = 31
    81116 = 44
    81100 = 181
    81119 = 22
    81118 = 47
    83038 = 1
    83037 = 2
    82891 = 6
    83504 = 28
    82824 = 68
    83507 = 24
    82489 = 1180
    81120 = 124
    81114 = 61
    82779 = 5895
    82940 = 953
    83142 = 1073
    82233 = 5491
    83271 = 66
    80528 = 133
    83508 = 77
    83301 = 401
}

function Update-Sheet([string]$sheetName) {
    $ws = $wb.Worksheets.Item($sheetName)

    # The row to drop is always the second row (first data row, right below
    # the header): it is the 2024-04-04 "Anime LIVE" entry in both affected
    # sheets. Deleting it shifts every following row up by one and Excel
    # automatically shrinks the sheet's dimension/used range.
    $ws.Rows.Item(2).Delete()

    # Figure out how many data rows remain after the delete.
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

    for ($r = 2; $r -le $lastRow; $r++) {
        # Renumber the running index in column A to match the new row
        # position (row 2 -> 1, row 3 -> 2, ...).
        $ws.Cells.Item($r, 1).Value = $r - 1

        # Refresh the "想去人数" count (column F) using the event's
        # bilibili show id, parsed out of the Link in column H.
        $link = $ws.Cells.Item($r, 8).Value()
        if ($link -match "id=(\d+)") {
            $id = [int]$matches[1]
            if ($idToNewF.ContainsKey($id)) {
                $ws.Cells.Item($r, 6).Value = $idToNewF[$id]
            }
        }
    }
}

Update-Sheet("展览")
Update-Sheet("全部类型")
